$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98

# Date-looking and number-looking text values need a quote prefix so Excel
# stores them as literal text instead of coercing to a date serial / number;
# ClearFormats() afterwards drops the resulting "quote prefix" cell style so
# the cell ends up with plain default formatting, same as the rest of the data.
$ws.Cells.Item($row, 1).Value = "'2024-01-27"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "23:20:50"
$ws.Cells.Item($row, 3).Value = "Saturday"

$ws.Cells.Item($row, 4).Value = "'03"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 137243
$ws.Cells.Item($row, 6).Value = 141764
$ws.Cells.Item($row, 7).Value = 171559
$ws.Cells.Item($row, 8).Value = 149119
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 121995
$ws.Cells.Item($row, 11).Value = 223875
$ws.Cells.Item($row, 12).Value = 257235
$ws.Cells.Item($row, 13).Value = 185462
$ws.Cells.Item($row, 14).Value = 110018
$ws.Cells.Item($row, 15).Value = 41430
$ws.Cells.Item($row, 16).Value = 30837
$ws.Cells.Item($row, 17).Value = 73627
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42489
$ws.Cells.Item($row, 20).Value = -1
